# Atualiza a planilha de notas: alunos que estavam "Aprovado" mas cuja
# media (P1, P2, P3) ficou abaixo do necessario passam para "Exame Final",
# recebendo a nota necessaria para aprovacao final na coluna H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha -> nova nota (coluna H) para os alunos que mudam de "Aprovado" para "Exame Final"
$updates = @{
    4  = 34
    5  = 9
    8  = 16
    10 = 26
    12 = 24
    13 = 35
    14 = 14
    15 = 12
    17 = 36
    19 = 40
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = "Exame Final"   # Coluna G - Situacao
    $ws.Cells.Item($row, 8).Value = $updates[$row]  # Coluna H - Nota para Aprovacao Final
}
